$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Staro"

# Duplicate the sheet to create "Novo", placed right after "Staro"
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Novo"

# --- Expand the DIZAJNIRANJE (design) section from 7 rows (7-13) to 15 rows (7-21) ---
# Insert 8 blank rows just above the old last row of that section (row 13), so that
# row keeps its "thick bottom border" formatting and ends up at row 21.
$ws2.Rows("13:20").Insert()

# Give the newly inserted rows (12-20) the same formatting as the regular design rows
# (row 12 currently still carries the old "second to last row" style, and rows 13-20
# are blank with generic styling) by cloning row 11's look-and-feel across them.
$ws2.Range("A12:C20").Value = 0
$ws2.Range("A11:D11").Copy()
$ws2.Range("A12:D20").PasteSpecial(-4122)
$ws2.Range("D12:D20").ClearContents()
for ($i = 12; $i -le 20; $i++) {
    $ws2.Range("D$i").Formula = "=B$i*C$i"
}

# --- Write the new DIZAJNIRANJE content (rows 7-20), then restore the final row ---
$ws2.Range("A7").Value = "Izrada modela podataka"
$ws2.Range("B7").Value = 1.5
$ws2.Range("C7").Value = 200

$ws2.Range("A8").Value = "Dizajn funkcije prijave u aplikaciju"
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = 200

$ws2.Range("A9").Value = "Dizajn funkcije upisa namirnica"
$ws2.Range("B9").Value = 2
$ws2.Range("C9").Value = 200

$ws2.Range("A10").Value = "Dizajn funkcije popisa namirnica"
$ws2.Range("B10").Value = 2
$ws2.Range("C10").Value = 200

$ws2.Range("A11").Value = "Dizajn funkcije spremanja u PDF"
$ws2.Range("B11").Value = 2
$ws2.Range("C11").Value = 200

$ws2.Range("A12").Value = "Dizajn funkcije rasporeda kupovine"
$ws2.Range("B12").Value = 2
$ws2.Range("C12").Value = 200

$ws2.Range("A13").Value = "Dizajn funkcije liste za kupovinu"
$ws2.Range("B13").Value = 2
$ws2.Range("C13").Value = 200

$ws2.Range("A14").Value = "Dizajn funkcije uređivanja liste za kupovinu"
$ws2.Range("B14").Value = 2
$ws2.Range("C14").Value = 200

$ws2.Range("A15").Value = "Dizajn funkcije predviđanja potrošnje"
$ws2.Range("B15").Value = 2
$ws2.Range("C15").Value = 200

$ws2.Range("A16").Value = "Dizajn funkcije odabira ispisa liste"
$ws2.Range("B16").Value = 2
$ws2.Range("C16").Value = 200

$ws2.Range("A17").Value = "Dizajn funkcije ispisa liste za kupovinu"
$ws2.Range("B17").Value = 2
$ws2.Range("C17").Value = 200

$ws2.Range("A18").Value = "Dizajn funkcije email podsjetnika"
$ws2.Range("B18").Value = 2
$ws2.Range("C18").Value = 200

$ws2.Range("A19").Value = "Dizajn funkcije domaćinstva"
$ws2.Range("B19").Value = 2
$ws2.Range("C19").Value = 200

$ws2.Range("A20").Value = "Dizajn funkcije prijave problema"
$ws2.Range("B20").Value = 2
$ws2.Range("C20").Value = 200

# Row 21 keeps "Izrada tehničke dokumentacije" (same wording as before), but BROJ SATI
# drops from 20 to 8.
$ws2.Range("A21").Value = "Izrada tehničke dokumentacije"
$ws2.Range("B21").Value = 8
$ws2.Range("C21").Value = 200

# --- Correct the TESTIRANJE unit prices (rows 41-43 after the row insert) ---
$ws2.Range("C41").Value = 200
$ws2.Range("C42").Value = 220
$ws2.Range("C43").Value = 290

# --- View/selection bookkeeping to mirror the authored workbook ---
$ws1.Range("A1").Select()
$ws2.Range("D48").Select()
$ws2.Activate()

$wb.Windows.Item(1).ActiveSheet.Range("D48").Select()
